$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-11-04 08:10:46"

foreach ($sheetName in @("Главные", "Линейные")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("AA2:AA26").Value = $newTimestamp
}
